# Refresh the recomputed aerodynamic comparison data (vortex-panel output
# vs. the Abbot & von Doenhoff reference) on Sheet1. Only the computed
# columns (B,C,D,E,G,H,J,K) change; the independent/reference columns
# (A = x, F = one_minus_Cp_Abbot, I = normalized_V_Abbott) stay as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 581.81184391588317
$ws.Range("D2").Value = 581.81184391588317
$ws.Range("E2").Value = -32.850502172080006
$ws.Range("G2").Value = 33.850502172080006
$ws.Range("J2").Value = 5.8181184391588321
$ws.Range("B3").Value = 52.468904602940491
$ws.Range("C3").Value = 62.015268114861158
$ws.Range("D3").Value = 81.23348711947915
$ws.Range("E3").Value = 0.3401120570409415
$ws.Range("G3").Value = 0.6598879429590585
$ws.Range("H3").Value = 0.031074910873528892
$ws.Range("J3").Value = 0.81233487119479153
$ws.Range("K3").Value = 0.015418588993489357
$ws.Range("B4").Value = 82.103732443261279
$ws.Range("C4").Value = 58.775745317600325
$ws.Range("D4").Value = 100.97331884589143
$ws.Range("E4").Value = -0.019561111875405102
$ws.Range("G4").Value = 1.0195611118754051
$ws.Range("H4").Value = 0.0094664474013911801
$ws.Range("J4").Value = 1.0097331884589142
$ws.Range("K4").Value = 0.00470964025762622
$ws.Range("B5").Value = 100.21699723365364
$ws.Range("C5").Value = 47.853232711935554
$ws.Range("D5").Value = 111.05574463085105
$ws.Range("E5").Value = -0.23333784155128012
$ws.Range("G5").Value = 1.2333378415512801
$ws.Range("H5").Value = 0.0061741808611764544
$ws.Range("J5").Value = 1.1105574463085106
$ws.Range("K5").Value = 0.003090263636884642
$ws.Range("B6").Value = 111.67353020854374
$ws.Range("C6").Value = 33.824354399947424
$ws.Range("D6").Value = 116.68360767396496
$ws.Range("E6").Value = -0.36150642998117721
$ws.Range("G6").Value = 1.3615064299811772
$ws.Range("H6").Value = 0.01196920901220804
$ws.Range("J6").Value = 1.1668360767396495
$ws.Range("K6").Value = 0.0061021492847959328
$ws.Range("B7").Value = 115.52112660204648
$ws.Range("C7").Value = 25.503532510051429
$ws.Range("D7").Value = 118.30283539246767
$ws.Range("E7").Value = -0.39955608618973026
$ws.Range("G7").Value = 1.3995560861897303
$ws.Range("H7").Value = 0.0017431624894933303
$ws.Range("J7").Value = 1.1830283539246766
$ws.Range("K7").Value = 0.00082064702307710802
$ws.Range("B8").Value = 117.16422706548609
$ws.Range("C8").Value = 19.770979568610642
$ws.Range("D8").Value = 118.82065366322136
$ws.Range("E8").Value = -0.41183477369552013
$ws.Range("G8").Value = 1.4118347736955201
$ws.Range("H8").Value = 0.00059161849434450514
$ws.Range("J8").Value = 1.1882065366322137
$ws.Range("K8").Value = 0.00017385238401827039
$ws.Range("B9").Value = 118.07231132749567
$ws.Range("C9").Value = 12.017606603073391
$ws.Range("D9").Value = 118.68232206475949
$ws.Range("E9").Value = -0.40854935706832962
$ws.Range("G9").Value = 1.4085493570683296
$ws.Range("H9").Value = 0.0017368128502270803
$ws.Range("J9").Value = 1.1868232206475948
$ws.Range("K9").Value = 0.00099055501044202176
$ws.Range("B10").Value = 117.71340178360596
$ws.Range("C10").Value = 6.7924605369065594
$ws.Range("D10").Value = 117.90921286996227
$ws.Range("E10").Value = -0.39025824796140784
$ws.Range("G10").Value = 1.3902582479614078
$ws.Range("H10").Value = 0.0062485718646119919
$ws.Range("J10").Value = 1.1790921286996228
$ws.Range("K10").Value = 0.0033033569741143581
$ws.Range("B11").Value = 116.82407772165713
$ws.Range("C11").Value = 2.9483476615003887
$ws.Range("D11").Value = 116.86127626142401
$ws.Range("E11").Value = -0.36565578894488637
$ws.Range("G11").Value = 1.3656557889448864
$ws.Range("H11").Value = 0.0089580631749735257
$ws.Range("J11").Value = 1.1686127626142402
$ws.Range("K11").Value = 0.0045887882331854743
$ws.Range("B12").Value = 115.67562454978497
$ws.Range("C12").Value = -0.0089964029785963243
$ws.Range("D12").Value = 115.67562489962214
$ws.Range("E12").Value = -0.33808501959180814
$ws.Range("G12").Value = 1.3380850195918081
$ws.Range("H12").Value = 0.008825911413475513
$ws.Range("J12").Value = 1.1567562489962213
$ws.Range("K12").Value = 0.0045126944955065466
$ws.Range("B13").Value = 113.05698777479228
$ws.Range("C13").Value = -4.2111426559548653
$ws.Range("D13").Value = 113.13538883646589
$ws.Range("E13").Value = -0.27996162071783304
$ws.Range("G13").Value = 1.279961620717833
$ws.Range("H13").Value = 0.0062409777035458041
$ws.Range("J13").Value = 1.131353888364659
$ws.Range("K13").Value = 0.0032124331588907571
$ws.Range("B14").Value = 110.33311265053909
$ws.Range("C14").Value = -6.963232862893924
$ws.Range("D14").Value = 110.55262257883997
$ws.Range("E14").Value = -0.22218823590594372
$ws.Range("G14").Value = 1.2221882359059437
$ws.Range("H14").Value = 0.0047327069169839251
$ws.Range("J14").Value = 1.1055262257883998
$ws.Range("K14").Value = 0.0022326482054154623
$ws.Range("B15").Value = 107.62656484123201
$ws.Range("C15").Value = -8.8419875975868418
$ws.Range("D15").Value = 107.98915780854945
$ws.Range("E15").Value = -0.16616582041997963
$ws.Range("G15").Value = 1.1661658204199796
$ws.Range("H15").Value = 0.00014221305315583266
$ws.Range("J15").Value = 1.0798915780854945
$ws.Range("K15").Value = 0.00010039066157921485
$ws.Range("B16").Value = 104.87491866295302
$ws.Range("C16").Value = -10.202988149107831
$ws.Range("D16").Value = 105.37005993987023
$ws.Range("E16").Value = -0.11028495317318465
$ws.Range("G16").Value = 1.1102849531731847
$ws.Range("H16").Value = 0.0011586593085524495
$ws.Range("J16").Value = 1.0537005993987023
$ws.Range("K16").Value = 0.00066533656097091446
$ws.Range("B17").Value = 101.80826909238152
$ws.Range("C17").Value = -11.276470728167068
$ws.Range("D17").Value = 102.43086667440618
$ws.Range("E17").Value = -0.049208244766997478
$ws.Range("G17").Value = 1.0492082447669975
$ws.Range("H17").Value = 0.0049887401982734083
$ws.Range("J17").Value = 1.0243086667440617
$ws.Range("K17").Value = 0.0022589694168901331
$ws.Range("B18").Value = 97.632286567408542
$ws.Range("C18").Value = -12.155052156447296
$ws.Range("D18").Value = 98.386018688157804
$ws.Range("E18").Value = 0.032019132669346306
$ws.Range("G18").Value = 0.96798086733065369
$ws.Range("H18").Value = 0.012532288002775872
$ws.Range("J18").Value = 0.98386018688157806
$ws.Range("K18").Value = 0.005992011126357953
$ws.Range("B19").Value = 94.259165602612342
$ws.Range("C19").Value = -12.442818184847773
$ws.Range("D19").Value = 95.076884806360141
$ws.Range("E19").Value = 0.096038597551812432
$ws.Range("G19").Value = 0.90396140244818757
$ws.Range("H19").Value = 0.0022501076730821854
$ws.Range("J19").Value = 0.9507688480636014
$ws.Range("K19").Value = 0.0012932268239480627
$ws.Range("B20").Value = 69.849837469570758
$ws.Range("C20").Value = -45.513414731382063
$ws.Range("D20").Value = 83.369483116043341
$ws.Range("E20").Value = 0.30495292849637634
$ws.Range("G20").Value = 0.69504707150362366
$ws.Range("J20").Value = 0.83369483116043341
